$wb = $excel.ActiveWorkbook

# Add a new worksheet named "campos" at the end of the workbook
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "campos"

# Populate the new sheet with the field values (row 4 intentionally left blank)
$newSheet.Range("A1").Value = "fila 1"
$newSheet.Range("A2").Value = "fila 2"
$newSheet.Range("A3").Value = "fila 3 "
$newSheet.Range("A5").Value = "fila 5"

# Update selections on each sheet to match the saved view state
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A2").Select() | Out-Null

$hoja2 = $wb.Worksheets.Item("hoja2")
$hoja2.Range("A2").Select() | Out-Null

# Make the new "campos" sheet the active sheet with A6 selected
$newSheet.Activate()
$newSheet.Range("A6").Select() | Out-Null
